$wb = $excel.ActiveWorkbook

# --- Ninja sheet: fill in a new assessment column (C) for every student ---
$wsNinja = $wb.Worksheets.Item("Ninja")

$wsNinja.Range("C2").Value = 1
$wsNinja.Range("C3").Value = 1
$wsNinja.Range("C4").Value = 1
$wsNinja.Range("C5").Value = 1
$wsNinja.Range("C6").Value = 0
$wsNinja.Range("C7").Value = 1
$wsNinja.Range("C8").Value = 1
$wsNinja.Range("C9").Value = 1
$wsNinja.Range("C10").Value = 1
$wsNinja.Range("C11").Value = 1
$wsNinja.Range("C12").Value = 1
$wsNinja.Range("C13").Value = 1
$wsNinja.Range("C14").Value = 1
$wsNinja.Range("C15").Value = 1
$wsNinja.Range("C16").Value = 1
$wsNinja.Range("C17").Value = 1
$wsNinja.Range("C18").Value = 1
$wsNinja.Range("C19").Value = 1
$wsNinja.Range("C20").Value = 0
$wsNinja.Range("C21").Value = 1

# --- Update the selection on each sheet; whichever sheet is activated last
#     becomes the workbook's active tab on save. ---

$wsSenador = $wb.Worksheets.Item("Senador")
$wsSenador.Activate() | Out-Null
$wsSenador.Range("C2").Select() | Out-Null

$wsNinja.Activate() | Out-Null
$wsNinja.Range("C1").Select() | Out-Null

$wsAstronauta = $wb.Worksheets.Item("Astronauta")
$wsAstronauta.Activate() | Out-Null
$wsAstronauta.Range("B24").Select() | Out-Null
